# Update out-of-date "UndoRedoStack" section of the Logic/Model class
# diagram to reflect the move to VersionedAddressBook, and refresh the
# auto-date footer fields that PowerPoint re-stamps whenever the deck is
# opened and saved.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Refresh the "datetimeFigureOut" footer date placeholders across the
#    slide master, every slide layout, and the notes master.
# ---------------------------------------------------------------------
$newDate = "4/16/2018"

function Update-DateText($shape) {
    if ($shape.HasTextFrame -eq -1 -and $shape.TextFrame.HasText -eq -1) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq "7/20/17") {
            $tr.Text = $newDate
        }
    }
}

$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    Update-DateText($master.Shapes.Item($i))
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        Update-DateText($layout.Shapes.Item($i))
    }
}

$notesMaster = $p.NotesMaster
$notesMaster.HeadersFooters.DateAndTime.Text = $newDate

# ---------------------------------------------------------------------
# 2. Remove the now out-of-date "UndoRedoStack" shapes from the
#    Logic/Model component diagram: the blue "UndoRedo / Stack"
#    rectangle, its connector arrow, and the "1" multiplicity label.
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Id -eq 59 -or $shp.Id -eq 61 -or $shp.Id -eq 63) {
        $shp.Delete()
    }
}
